# Update the dSF (column F) values on Sheet1 with re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 2
$ws.Range("F22").Value = -1
$ws.Range("F23").Value = 8
$ws.Range("F26").Value = -2
$ws.Range("F27").Value = -3
$ws.Range("F28").Value = -2
$ws.Range("F29").Value = 2
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = -2
$ws.Range("F32").Value = -2
$ws.Range("F34").Value = 2
$ws.Range("F35").Value = -2
$ws.Range("F36").Value = 1
$ws.Range("F37").Value = -2
$ws.Range("F38").Value = 3
